$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header string stays the same text; rewrite it so the workbook keeps
# it as a (shared) string cell.
$ws.Range("A1").Value = "HK_G_acc_G"

# Updated accuracy values (A2:A49) reflecting the new results.
$ws.Range("A2").Value = 51.179673321234119
$ws.Range("A3").Value = 51.179673321234119
$ws.Range("A4").Value = 51.179673321234119
$ws.Range("A5").Value = 51.361161524500908
$ws.Range("A6").Value = 51.542649727767696
$ws.Range("A7").Value = 51.542649727767696
$ws.Range("A8").Value = 52.813067150635206
$ws.Range("A9").Value = 52.631578947368418
$ws.Range("A10").Value = 52.45009074410163
$ws.Range("A11").Value = 52.45009074410163
$ws.Range("A12").Value = 52.08711433756806
$ws.Range("A13").Value = 52.994555353901994
$ws.Range("A14").Value = 52.631578947368418
$ws.Range("A15").Value = 52.631578947368418
$ws.Range("A16").Value = 53.176043557168782
$ws.Range("A17").Value = 52.08711433756806
$ws.Range("A18").Value = 52.994555353901994
$ws.Range("A19").Value = 52.268602540834841
$ws.Range("A20").Value = 52.08711433756806
$ws.Range("A21").Value = 52.08711433756806
$ws.Range("A22").Value = 52.268602540834841
$ws.Range("A23").Value = 50.453720508166967
$ws.Range("A24").Value = 49.364791288566245
$ws.Range("A25").Value = 49.727767695099814
$ws.Range("A26").Value = 51.905626134301272
$ws.Range("A27").Value = 51.905626134301272
$ws.Range("A28").Value = 52.631578947368418
$ws.Range("A29").Value = 53.539019963702358
$ws.Range("A30").Value = 52.813067150635206
$ws.Range("A31").Value = 53.176043557168782
$ws.Range("A32").Value = 50.090744101633391
$ws.Range("A33").Value = 50.635208711433755
$ws.Range("A34").Value = 50.816696914700543
$ws.Range("A35").Value = 51.542649727767696
$ws.Range("A36").Value = 51.542649727767696
$ws.Range("A37").Value = 56.805807622504531
$ws.Range("A38").Value = 50.090744101633391
$ws.Range("A39").Value = 51.179673321234119
$ws.Range("A40").Value = 51.179673321234119
$ws.Range("A41").Value = 52.631578947368418
$ws.Range("A42").Value = 52.45009074410163
$ws.Range("A43").Value = 53.35753176043557
$ws.Range("A44").Value = 52.994555353901994
$ws.Range("A45").Value = 51.361161524500908
$ws.Range("A46").Value = 51.361161524500908
$ws.Range("A47").Value = 50.998185117967331
$ws.Range("A48").Value = 54.446460980036292
$ws.Range("A49").Value = 52.08711433756806
